$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 6800016
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "TNS"
$arr[0,1] = "Connahs Quay"
$arr[0,2] = 6
$arr[0,3] = 2
$arr[0,4] = "H"
$arr[0,5] = 1.25
$arr[0,6] = 5
$arr[0,7] = 9
$arr[0,8] = 1.222
$arr[0,9] = 5
$arr[0,10] = 11
$arr[0,11] = -1.75
$arr[0,12] = 1.925
$arr[0,13] = 1.875
$arr[0,14] = 2.75
$arr[0,15] = 1.775
$arr[0,16] = 2.025
$arr[0,17] = 0.222
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.925
$arr[0,21] = -1
$arr[0,22] = 0.7749999999999999
$arr[0,23] = -1
$ws.Range("E2:AB2").Value2 = $arr

$ws.Range("B3").Value2 = 6800423
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Haverfordwest County"
$arr[0,1] = "Pontypridd Town"
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = "D"
$arr[0,5] = 1.75
$arr[0,6] = 3.6
$arr[0,7] = 3.75
$arr[0,8] = 2.05
$arr[0,9] = 3.5
$arr[0,10] = 2.9
$arr[0,11] = -0.25
$arr[0,12] = 1.8
$arr[0,13] = 2
$arr[0,14] = 2.5
$arr[0,15] = 1.975
$arr[0,16] = 1.825
$arr[0,17] = -1
$arr[0,18] = 2.5
$arr[0,19] = -1
$arr[0,20] = -0.5
$arr[0,21] = 0.5
$arr[0,22] = -1
$arr[0,23] = 0.825
$ws.Range("E3:AB3").Value2 = $arr

$ws.Range("B26").Value2 = 6800443
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Newtown"
$arr[0,1] = "Cardiff MU"
$arr[0,2] = 2
$arr[0,3] = 1
$arr[0,4] = "H"
$arr[0,5] = 2.7
$arr[0,6] = 3.25
$arr[0,7] = 2.3
$arr[0,8] = 2.4
$arr[0,9] = 3.25
$arr[0,10] = 2.7
$arr[0,11] = 0
$arr[0,12] = 1.775
$arr[0,13] = 2.025
$arr[0,14] = 2.5
$arr[0,15] = 1.95
$arr[0,16] = 1.85
$arr[0,17] = 1.4
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.7749999999999999
$arr[0,21] = -1
$arr[0,22] = 0.95
$arr[0,23] = -1
$ws.Range("E26:AB26").Value2 = $arr

$ws.Range("B27").Value2 = 6800022
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Caernarfon Town"
$arr[0,1] = "Connahs Quay"
$arr[0,2] = 0
$arr[0,3] = 4
$arr[0,4] = "A"
$arr[0,5] = 3.5
$arr[0,6] = 3.6
$arr[0,7] = 1.8
$arr[0,8] = 3.5
$arr[0,9] = 3.5
$arr[0,10] = 1.909
$arr[0,11] = 0.5
$arr[0,12] = 1.85
$arr[0,13] = 1.95
$arr[0,14] = 2.75
$arr[0,15] = 1.8
$arr[0,16] = 2
$arr[0,17] = -1
$arr[0,18] = -1
$arr[0,19] = 0.909
$arr[0,20] = -1
$arr[0,21] = 0.95
$arr[0,22] = 0.8
$arr[0,23] = -1
$ws.Range("E27:AB27").Value2 = $arr

$ws.Range("B28").Value2 = 6800023
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "TNS"
$arr[0,1] = "Aberystwyth"
$arr[0,2] = 6
$arr[0,3] = 0
$arr[0,4] = "H"
$arr[0,5] = 1.02
$arr[0,6] = 11
$arr[0,7] = 26
$arr[0,8] = 1.02
$arr[0,9] = 21
$arr[0,10] = 41
$arr[0,11] = -4.25
$arr[0,12] = 1.825
$arr[0,13] = 1.975
$arr[0,14] = 5.25
$arr[0,15] = 1.975
$arr[0,16] = 1.825
$arr[0,17] = 0.02000000000000002
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.825
$arr[0,21] = -1
$arr[0,22] = 0.9750000000000001
$arr[0,23] = -1
$ws.Range("E28:AB28").Value2 = $arr

$ws.Range("B58").Value2 = 6800464
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Colwyn Bay"
$arr[0,1] = "Cardiff MU"
$arr[0,2] = 2
$arr[0,3] = 2
$arr[0,4] = "D"
$arr[0,5] = 3
$arr[0,6] = 3.3
$arr[0,7] = 2.1
$arr[0,8] = 3.8
$arr[0,9] = 3.6
$arr[0,10] = 1.727
$arr[0,11] = 0.5
$arr[0,12] = 2.025
$arr[0,13] = 1.775
$arr[0,14] = 2.75
$arr[0,15] = 2
$arr[0,16] = 1.8
$arr[0,17] = -1
$arr[0,18] = 2.6
$arr[0,19] = -1
$arr[0,20] = 1.025
$arr[0,21] = -1
$arr[0,22] = 1
$arr[0,23] = -1
$ws.Range("E58:AB58").Value2 = $arr

$ws.Range("B59").Value2 = 6800466
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Haverfordwest County"
$arr[0,1] = "Aberystwyth"
$arr[0,2] = 3
$arr[0,3] = 0
$arr[0,4] = "H"
$arr[0,5] = 1.7
$arr[0,6] = 3.6
$arr[0,7] = 4.1
$arr[0,8] = 1.55
$arr[0,9] = 3.6
$arr[0,10] = 5.25
$arr[0,11] = -1
$arr[0,12] = 2.025
$arr[0,13] = 1.775
$arr[0,14] = 2.5
$arr[0,15] = 1.9
$arr[0,16] = 1.9
$arr[0,17] = 0.55
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 1.025
$arr[0,21] = -1
$arr[0,22] = 0.8999999999999999
$arr[0,23] = -1
$ws.Range("E59:AB59").Value2 = $arr

$ws.Range("B60").Value2 = 6800032
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "TNS"
$arr[0,1] = "Pontypridd Town"
$arr[0,2] = 4
$arr[0,3] = 0
$arr[0,4] = "H"
$arr[0,5] = 1.05
$arr[0,6] = 15
$arr[0,7] = 29
$arr[0,8] = 1.055
$arr[0,9] = 12
$arr[0,10] = 34
$arr[0,11] = -3.25
$arr[0,12] = 2
$arr[0,13] = 1.8
$arr[0,14] = 4
$arr[0,15] = 1.925
$arr[0,16] = 1.875
$arr[0,17] = 0.05499999999999994
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 1
$arr[0,21] = -1
$arr[0,22] = 0
$arr[0,23] = 0
$ws.Range("E60:AB60").Value2 = $arr

$ws.Range("B68").Value2 = 6800475
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Haverfordwest County"
$arr[0,1] = "Penybont"
$arr[0,2] = 3
$arr[0,3] = 2
$arr[0,4] = "H"
$arr[0,5] = 2.625
$arr[0,6] = 3.2
$arr[0,7] = 2.375
$arr[0,8] = 3
$arr[0,9] = 3
$arr[0,10] = 2.15
$arr[0,11] = 0.25
$arr[0,12] = 1.9
$arr[0,13] = 1.9
$arr[0,14] = 2.25
$arr[0,15] = 1.85
$arr[0,16] = 1.95
$arr[0,17] = 2
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.8999999999999999
$arr[0,21] = -1
$arr[0,22] = 0.8500000000000001
$arr[0,23] = -1
$ws.Range("E68:AB68").Value2 = $arr

$ws.Range("B69").Value2 = 6800036
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Caernarfon Town"
$arr[0,1] = "Aberystwyth"
$arr[0,2] = 3
$arr[0,3] = 0
$arr[0,4] = "H"
$arr[0,5] = 1.4
$arr[0,6] = 4.5
$arr[0,7] = 5.75
$arr[0,8] = 1.363
$arr[0,9] = 4.75
$arr[0,10] = 6
$arr[0,11] = -1.25
$arr[0,12] = 1.875
$arr[0,13] = 1.925
$arr[0,14] = 3
$arr[0,15] = 1.775
$arr[0,16] = 2.025
$arr[0,17] = 0.363
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.875
$arr[0,21] = -1
$arr[0,22] = 0
$arr[0,23] = 0
$ws.Range("E69:AB69").Value2 = $arr

$ws.Range("B70").Value2 = 6800472
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Barry Town"
$arr[0,1] = "Pontypridd Town"
$arr[0,2] = 2
$arr[0,3] = 0
$arr[0,4] = "H"
$arr[0,5] = 2.6
$arr[0,6] = 3.2
$arr[0,7] = 2.5
$arr[0,8] = 2.6
$arr[0,9] = 3.2
$arr[0,10] = 2.5
$arr[0,11] = 0
$arr[0,12] = 1.975
$arr[0,13] = 1.825
$arr[0,14] = 2.5
$arr[0,15] = 1.95
$arr[0,16] = 1.85
$arr[0,17] = 1.6
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.9750000000000001
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.8500000000000001
$ws.Range("E70:AB70").Value2 = $arr

$ws.Range("B86").Value2 = 6800484
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Bala Town"
$arr[0,1] = "Pontypridd Town"
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = "D"
$arr[0,5] = 1.615
$arr[0,6] = 3.5
$arr[0,7] = 4.75
$arr[0,8] = 1.533
$arr[0,9] = 3.75
$arr[0,10] = 6.5
$arr[0,11] = -1
$arr[0,12] = 1.95
$arr[0,13] = 1.85
$arr[0,14] = 2.25
$arr[0,15] = 1.775
$arr[0,16] = 2.025
$arr[0,17] = -1
$arr[0,18] = 2.75
$arr[0,19] = -1
$arr[0,20] = -1
$arr[0,21] = 0.8500000000000001
$arr[0,22] = -1
$arr[0,23] = 1.025
$ws.Range("E86:AB86").Value2 = $arr

$ws.Range("B87").Value2 = 6800042
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Newtown"
$arr[0,1] = "TNS"
$arr[0,2] = 0
$arr[0,3] = 2
$arr[0,4] = "A"
$arr[0,5] = 7
$arr[0,6] = 6
$arr[0,7] = 1.25
$arr[0,8] = 6.5
$arr[0,9] = 5.75
$arr[0,10] = 1.3
$arr[0,11] = 1.5
$arr[0,12] = 1.975
$arr[0,13] = 1.825
$arr[0,14] = 3.25
$arr[0,15] = 1.85
$arr[0,16] = 1.95
$arr[0,17] = -1
$arr[0,18] = -1
$arr[0,19] = 0.3
$arr[0,20] = -1
$arr[0,21] = 0.825
$arr[0,22] = -1
$arr[0,23] = 0.95
$ws.Range("E87:AB87").Value2 = $arr

$ws.Range("B97").Value2 = 6800493
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Newtown"
$arr[0,1] = "Pontypridd Town"
$arr[0,2] = 3
$arr[0,3] = 1
$arr[0,4] = "H"
$arr[0,5] = 1.45
$arr[0,6] = 3.75
$arr[0,7] = 6.5
$arr[0,8] = 1.444
$arr[0,9] = 3.8
$arr[0,10] = 7.5
$arr[0,11] = -1.25
$arr[0,12] = 2.025
$arr[0,13] = 1.775
$arr[0,14] = 2.5
$arr[0,15] = 2
$arr[0,16] = 1.8
$arr[0,17] = 0.444
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 1.025
$arr[0,21] = -1
$arr[0,22] = 1
$arr[0,23] = -1
$ws.Range("E97:AB97").Value2 = $arr

$ws.Range("B98").Value2 = 6800492
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Barry Town"
$arr[0,1] = "Penybont"
$arr[0,2] = 1
$arr[0,3] = 1
$arr[0,4] = "D"
$arr[0,5] = 3.75
$arr[0,6] = 3.75
$arr[0,7] = 1.727
$arr[0,8] = 3
$arr[0,9] = 3.6
$arr[0,10] = 2.05
$arr[0,11] = 0.25
$arr[0,12] = 2
$arr[0,13] = 1.8
$arr[0,14] = 2.75
$arr[0,15] = 1.875
$arr[0,16] = 1.925
$arr[0,17] = -1
$arr[0,18] = 2.6
$arr[0,19] = -1
$arr[0,20] = 0.5
$arr[0,21] = -0.5
$arr[0,22] = -1
$arr[0,23] = 0.925
$ws.Range("E98:AB98").Value2 = $arr

$ws.Range("B99").Value2 = 6800046
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Bala Town"
$arr[0,1] = "Caernarfon Town"
$arr[0,2] = 1
$arr[0,3] = 1
$arr[0,4] = "D"
$arr[0,5] = 2
$arr[0,6] = 3.6
$arr[0,7] = 3
$arr[0,8] = 1.727
$arr[0,9] = 3.8
$arr[0,10] = 3.8
$arr[0,11] = -0.75
$arr[0,12] = 1.95
$arr[0,13] = 1.85
$arr[0,14] = 3
$arr[0,15] = 1.925
$arr[0,16] = 1.875
$arr[0,17] = -1
$arr[0,18] = 2.8
$arr[0,19] = -1
$arr[0,20] = -1
$arr[0,21] = 0.8500000000000001
$arr[0,22] = -1
$arr[0,23] = 0.875
$ws.Range("E99:AB99").Value2 = $arr

$ws.Range("B115").Value2 = 6800507
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Newtown"
$arr[0,1] = "Connahs Quay"
$arr[0,2] = 1
$arr[0,3] = 2
$arr[0,4] = "A"
$arr[0,5] = 3.3
$arr[0,6] = 3.4
$arr[0,7] = 2
$arr[0,8] = 3.6
$arr[0,9] = 3.4
$arr[0,10] = 1.909
$arr[0,11] = 0.5
$arr[0,12] = 1.9
$arr[0,13] = 1.9
$arr[0,14] = 2.75
$arr[0,15] = 2
$arr[0,16] = 1.8
$arr[0,17] = -1
$arr[0,18] = -1
$arr[0,19] = 0.909
$arr[0,20] = -1
$arr[0,21] = 0.8999999999999999
$arr[0,22] = 0.5
$arr[0,23] = -0.5
$ws.Range("E115:AB115").Value2 = $arr

$ws.Range("B116").Value2 = 6800512
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Pontypridd Town"
$arr[0,1] = "Aberystwyth"
$arr[0,2] = 2
$arr[0,3] = 0
$arr[0,4] = "H"
$arr[0,5] = 2.4
$arr[0,6] = 4
$arr[0,7] = 2.25
$arr[0,8] = 1.909
$arr[0,9] = 3.6
$arr[0,10] = 3.6
$arr[0,11] = -0.5
$arr[0,12] = 1.925
$arr[0,13] = 1.875
$arr[0,14] = 2.25
$arr[0,15] = 1.85
$arr[0,16] = 1.95
$arr[0,17] = 0.909
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.925
$arr[0,21] = -1
$arr[0,22] = -0.5
$arr[0,23] = 0.475
$ws.Range("E116:AB116").Value2 = $arr

$ws.Range("B137").Value2 = 7721604
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Colwyn Bay"
$arr[0,1] = "Penybont"
$arr[0,2] = 1
$arr[0,3] = 2
$arr[0,4] = "A"
$arr[0,5] = 4.2
$arr[0,6] = 3.6
$arr[0,7] = 1.7
$arr[0,8] = 4
$arr[0,9] = 3.5
$arr[0,10] = 1.75
$arr[0,11] = 0.5
$arr[0,12] = 2
$arr[0,13] = 1.8
$arr[0,14] = 2.75
$arr[0,15] = 1.9
$arr[0,16] = 1.9
$arr[0,17] = -1
$arr[0,18] = -1
$arr[0,19] = 0.75
$arr[0,20] = -1
$arr[0,21] = 0.8
$arr[0,22] = 0.45
$arr[0,23] = -0.5
$ws.Range("E137:AB137").Value2 = $arr

$ws.Range("B138").Value2 = 7721563
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Bala Town"
$arr[0,1] = "Cardiff MU"
$arr[0,2] = 1
$arr[0,3] = 1
$arr[0,4] = "D"
$arr[0,5] = 1.833
$arr[0,6] = 3.4
$arr[0,7] = 3.75
$arr[0,8] = 1.666
$arr[0,9] = 3.5
$arr[0,10] = 4.5
$arr[0,11] = -0.75
$arr[0,12] = 1.9
$arr[0,13] = 1.9
$arr[0,14] = 2.25
$arr[0,15] = 1.825
$arr[0,16] = 1.975
$arr[0,17] = -1
$arr[0,18] = 2.5
$arr[0,19] = -1
$arr[0,20] = -1
$arr[0,21] = 0.8999999999999999
$arr[0,22] = -0.5
$arr[0,23] = 0.4875
$ws.Range("E138:AB138").Value2 = $arr

$ws.Range("B152").Value2 = 7721589
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Caernarfon Town"
$arr[0,1] = "Newtown"
$arr[0,2] = 1
$arr[0,3] = 0
$arr[0,4] = "H"
$arr[0,5] = 2.35
$arr[0,6] = 3.6
$arr[0,7] = 2.45
$arr[0,8] = 2.375
$arr[0,9] = 3.8
$arr[0,10] = 2.4
$arr[0,11] = 0
$arr[0,12] = 1.9
$arr[0,13] = 1.9
$arr[0,14] = 3.25
$arr[0,15] = 1.825
$arr[0,16] = 1.975
$arr[0,17] = 1.375
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.8999999999999999
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.9750000000000001
$ws.Range("E152:AB152").Value2 = $arr

$ws.Range("B153").Value2 = 7721588
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Bala Town"
$arr[0,1] = "Connahs Quay"
$arr[0,2] = 1
$arr[0,3] = 0
$arr[0,4] = "H"
$arr[0,5] = 3.75
$arr[0,6] = 3.8
$arr[0,7] = 1.7
$arr[0,8] = 2.55
$arr[0,9] = 3.2
$arr[0,10] = 2.45
$arr[0,11] = 0
$arr[0,12] = 1.9
$arr[0,13] = 1.9
$arr[0,14] = 2.5
$arr[0,15] = 1.875
$arr[0,16] = 1.925
$arr[0,17] = 1.55
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.8999999999999999
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.925
$ws.Range("E153:AB153").Value2 = $arr

$ws.Range("B163").Value2 = 7721617
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Haverfordwest County"
$arr[0,1] = "Aberystwyth"
$arr[0,2] = 2
$arr[0,3] = 0
$arr[0,4] = "H"
$arr[0,5] = 1.5
$arr[0,6] = 4.2
$arr[0,7] = 5.5
$arr[0,8] = 2.1
$arr[0,9] = 3.2
$arr[0,10] = 3.2
$arr[0,11] = -0.25
$arr[0,12] = 1.825
$arr[0,13] = 1.975
$arr[0,14] = 2.25
$arr[0,15] = 2
$arr[0,16] = 1.8
$arr[0,17] = 1.1
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.825
$arr[0,21] = -1
$arr[0,22] = -0.5
$arr[0,23] = 0.4
$ws.Range("E163:AB163").Value2 = $arr

$ws.Range("B164").Value2 = 7721592
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Caernarfon Town"
$arr[0,1] = "Connahs Quay"
$arr[0,2] = 2
$arr[0,3] = 1
$arr[0,4] = "H"
$arr[0,5] = 3.6
$arr[0,6] = 3.6
$arr[0,7] = 1.8
$arr[0,8] = 3.8
$arr[0,9] = 4
$arr[0,10] = 1.666
$arr[0,11] = 0.75
$arr[0,12] = 1.925
$arr[0,13] = 1.875
$arr[0,14] = 3.25
$arr[0,15] = 2
$arr[0,16] = 1.8
$arr[0,17] = 2.8
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.925
$arr[0,21] = -1
$arr[0,22] = -0.5
$arr[0,23] = 0.4
$ws.Range("E164:AB164").Value2 = $arr

$ws.Range("B170").Value2 = 7721594
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Connahs Quay"
$arr[0,1] = "Newtown"
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = "D"
$arr[0,5] = 1.4
$arr[0,6] = 4.5
$arr[0,7] = 5.75
$arr[0,8] = 1.55
$arr[0,9] = 4.5
$arr[0,10] = 4.5
$arr[0,11] = -1
$arr[0,12] = 1.9
$arr[0,13] = 1.9
$arr[0,14] = 3
$arr[0,15] = 1.8
$arr[0,16] = 2
$arr[0,17] = -1
$arr[0,18] = 3.5
$arr[0,19] = -1
$arr[0,20] = -1
$arr[0,21] = 0.8999999999999999
$arr[0,22] = -1
$arr[0,23] = 1
$ws.Range("E170:AB170").Value2 = $arr

$ws.Range("B171").Value2 = 7721620
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Barry Town"
$arr[0,1] = "Penybont"
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = "D"
$arr[0,5] = 3.1
$arr[0,6] = 3.1
$arr[0,7] = 2.15
$arr[0,8] = 4
$arr[0,9] = 3.2
$arr[0,10] = 1.909
$arr[0,11] = 0.5
$arr[0,12] = 1.875
$arr[0,13] = 1.925
$arr[0,14] = 2.5
$arr[0,15] = 1.975
$arr[0,16] = 1.825
$arr[0,17] = -1
$arr[0,18] = 2.2
$arr[0,19] = -1
$arr[0,20] = 0.875
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.825
$ws.Range("E171:AB171").Value2 = $arr

$ws.Range("B174").Value2 = 7721596
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Bala Town"
$arr[0,1] = "Newtown"
$arr[0,2] = 1
$arr[0,3] = 1
$arr[0,4] = "D"
$arr[0,5] = 2
$arr[0,6] = 3.5
$arr[0,7] = 3
$arr[0,8] = 2
$arr[0,9] = 3.4
$arr[0,10] = 3
$arr[0,11] = -0.25
$arr[0,12] = 1.85
$arr[0,13] = 1.95
$arr[0,14] = 2.5
$arr[0,15] = 1.925
$arr[0,16] = 1.875
$arr[0,17] = -1
$arr[0,18] = 2.4
$arr[0,19] = -1
$arr[0,20] = -0.5
$arr[0,21] = 0.475
$arr[0,22] = -1
$arr[0,23] = 0.875
$ws.Range("E174:AB174").Value2 = $arr

$ws.Range("B175").Value2 = 7721597
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "TNS"
$arr[0,1] = "Caernarfon Town"
$arr[0,2] = 7
$arr[0,3] = 1
$arr[0,4] = "H"
$arr[0,5] = 1.1
$arr[0,6] = 9
$arr[0,7] = 11
$arr[0,8] = 1.1
$arr[0,9] = 10
$arr[0,10] = 17
$arr[0,11] = -2.75
$arr[0,12] = 1.975
$arr[0,13] = 1.825
$arr[0,14] = 4
$arr[0,15] = 1.95
$arr[0,16] = 1.85
$arr[0,17] = 0.1000000000000001
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.9750000000000001
$arr[0,21] = -1
$arr[0,22] = 0.95
$arr[0,23] = -1
$ws.Range("E175:AB175").Value2 = $arr

$ws.Range("B177").Value2 = 7721625
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Colwyn Bay"
$arr[0,1] = "Pontypridd Town"
$arr[0,2] = 1
$arr[0,3] = 0
$arr[0,4] = "H"
$arr[0,5] = 2.45
$arr[0,6] = 3.4
$arr[0,7] = 2.45
$arr[0,8] = 2.625
$arr[0,9] = 3.4
$arr[0,10] = 2.3
$arr[0,11] = 0.25
$arr[0,12] = 1.8
$arr[0,13] = 2
$arr[0,14] = 2.5
$arr[0,15] = 1.975
$arr[0,16] = 1.825
$arr[0,17] = 1.625
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.8
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.825
$ws.Range("E177:AB177").Value2 = $arr

$ws.Range("B178").Value2 = 7721570
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Cardiff MU"
$arr[0,1] = "Connahs Quay"
$arr[0,2] = 0
$arr[0,3] = 3
$arr[0,4] = "A"
$arr[0,5] = 3.5
$arr[0,6] = 3.5
$arr[0,7] = 1.833
$arr[0,8] = 3.1
$arr[0,9] = 3.4
$arr[0,10] = 2
$arr[0,11] = 0.25
$arr[0,12] = 1.975
$arr[0,13] = 1.825
$arr[0,14] = 2.5
$arr[0,15] = 1.8
$arr[0,16] = 2
$arr[0,17] = -1
$arr[0,18] = -1
$arr[0,19] = 1
$arr[0,20] = -1
$arr[0,21] = 0.825
$arr[0,22] = 0.8
$arr[0,23] = -1
$ws.Range("E178:AB178").Value2 = $arr

$ws.Range("B188").Value2 = 7721572
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Cardiff MU"
$arr[0,1] = "Newtown"
$arr[0,2] = 1
$arr[0,3] = 3
$arr[0,4] = "A"
$arr[0,5] = 2.4
$arr[0,6] = 3.6
$arr[0,7] = 2.4
$arr[0,8] = 3.3
$arr[0,9] = 3.8
$arr[0,10] = 1.8
$arr[0,11] = 0.5
$arr[0,12] = 1.95
$arr[0,13] = 1.85
$arr[0,14] = 3
$arr[0,15] = 1.875
$arr[0,16] = 1.925
$arr[0,17] = -1
$arr[0,18] = -1
$arr[0,19] = 0.8
$arr[0,20] = -1
$arr[0,21] = 0.8500000000000001
$arr[0,22] = 0.875
$arr[0,23] = -1
$ws.Range("E188:AB188").Value2 = $arr

$ws.Range("B189").Value2 = 7721631
$arr = New-Object 'object[,]' 1,24
$arr[0,0] = "Penybont"
$arr[0,1] = "Haverfordwest County"
$arr[0,2] = 1
$arr[0,3] = 0
$arr[0,4] = "H"
$arr[0,5] = 1.8
$arr[0,6] = 3.6
$arr[0,7] = 3.6
$arr[0,8] = 1.7
$arr[0,9] = 3.75
$arr[0,10] = 4
$arr[0,11] = -0.75
$arr[0,12] = 1.95
$arr[0,13] = 1.85
$arr[0,14] = 2.5
$arr[0,15] = 1.925
$arr[0,16] = 1.875
$arr[0,17] = 0.7
$arr[0,18] = -1
$arr[0,19] = -1
$arr[0,20] = 0.475
$arr[0,21] = -0.5
$arr[0,22] = -1
$arr[0,23] = 0.875
$ws.Range("E189:AB189").Value2 = $arr
